$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do- FY16 Release")
$ws.Select()
